$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$t.Cell(2,1).Range.Text = "Lando Norris"
$t.Cell(2,4).Range.Text = "2"
$t.Cell(2,5).Range.Text = "0"
$t.Cell(2,6).Range.Text = "232.594"
$t.Cell(2,7).Range.Text = "-1.891"
$t.Cell(3,1).Range.Text = "Max Verstappen"
$t.Cell(3,4).Range.Text = "1"
$t.Cell(3,5).Range.Text = "-1"
$t.Cell(3,6).Range.Text = "235.359"
$t.Cell(3,7).Range.Text = "0.874"
$t.Cell(4,5).Range.Text = "2"
$t.Cell(4,6).Range.Text = "225.508"
$t.Cell(4,7).Range.Text = "-8.977"
$t.Cell(5,6).Range.Text = "230.886"
$t.Cell(5,7).Range.Text = "-3.599"
$t.Cell(6,5).Range.Text = "-3"
$t.Cell(6,6).Range.Text = "229.983"
$t.Cell(6,7).Range.Text = "-4.502"
$t.Cell(7,6).Range.Text = "239.019"
$t.Cell(7,7).Range.Text = "4.534"
$t.Cell(8,5).Range.Text = "0"
$t.Cell(8,6).Range.Text = "231.783"
$t.Cell(8,7).Range.Text = "-2.702"
$t.Cell(9,1).Range.Text = "Pierre Gasly"
$t.Cell(9,4).Range.Text = "9"
$t.Cell(9,5).Range.Text = "8"
$t.Cell(9,6).Range.Text = "218.210"
$t.Cell(9,7).Range.Text = "-16.275"
$t.Cell(10,1).Range.Text = "Sergio Perez"
$t.Cell(10,4).Range.Text = "8"
$t.Cell(10,5).Range.Text = "-1"
$t.Cell(10,6).Range.Text = "221.942"
$t.Cell(10,7).Range.Text = "-12.543"
$t.Cell(11,7).Range.Text = "-0.945"
$t.Cell(12,5).Range.Text = "-3"
$t.Cell(12,6).Range.Text = "241.081"
$t.Cell(12,7).Range.Text = "6.596"
$t.Cell(13,6).Range.Text = "313.419"
$t.Cell(13,7).Range.Text = "-0.840"
$t.Cell(14,5).Range.Text = "-4"
$t.Cell(14,6).Range.Text = "309.830"
$t.Cell(14,7).Range.Text = "-4.429"
$t.Cell(15,5).Range.Text = "-4"
$t.Cell(15,6).Range.Text = "283.702"
$t.Cell(15,7).Range.Text = "-30.557"
$t.Cell(16,5).Range.Text = "0"
$t.Cell(16,6).Range.Text = "307.581"
$t.Cell(16,7).Range.Text = "-6.679"
$t.Cell(17,1).Range.Text = "Guanyu Zhou"
$t.Cell(17,4).Range.Text = "13"
$t.Cell(17,5).Range.Text = "5"
$t.Cell(17,6).Range.Text = "339.687"
$t.Cell(17,7).Range.Text = "25.428"
$t.Cell(18,5).Range.Text = "2"
$t.Cell(18,6).Range.Text = "301.849"
$t.Cell(18,7).Range.Text = "-12.410"
$t.Cell(19,1).Range.Text = "Alexander Albon"
$t.Cell(19,4).Range.Text = "18"
$t.Cell(19,5).Range.Text = "-5"
$t.Cell(19,6).Range.Text = "303.891"
$t.Cell(19,7).Range.Text = "-10.368"
$t.Cell(20,5).Range.Text = "-3"
$t.Cell(20,6).Range.Text = "316.130"
$t.Cell(20,7).Range.Text = "1.871"
$t.Cell(21,5).Range.Text = "-4"
$t.Cell(21,6).Range.Text = "387.295"
$t.Cell(21,7).Range.Text = "-6.431"
